$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H135").Value = 3538.9473
$ws_ALC.Range("I135").Value = 3725.8823
$ws_ALC.Range("J135").Value = 1950
$ws_ALC.Range("K135").Value = 33532.9407
$ws_ALC.Range("L135").Value = 17550
$ws_ALC.Range("M135").Value = -30997.9407
$ws_ALC.Range("N135").Value = -22620

$ws_ALC.Range("H137").Value = 563424.25
$ws_ALC.Range("I137").Value = 957833.4
$ws_ALC.Range("J137").Value = 21111.75
$ws_ALC.Range("K137").Value = 2873500.2
$ws_ALC.Range("L137").Value = 63335.25
$ws_ALC.Range("M137").Value = -2870950.2
$ws_ALC.Range("N137").Value = -68435.25

$ws_ALC.Range("H138").Value = 5314.8447
$ws_ALC.Range("I138").Value = 0
$ws_ALC.Range("K138").Value = 0
$ws_ALC.Range("M138").Value = ""

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H2").Value = 13726.182
$ws_ARM.Range("J2").Value = 3443.111
$ws_ARM.Range("L2").Value = 3443.111
$ws_ARM.Range("N2").Value = -3669.111

$ws_ARM.Range("H32").Value = 2187.0352
$ws_ARM.Range("I32").Value = 2176.1072
$ws_ARM.Range("K32").Value = 2176.1072
$ws_ARM.Range("M32").Value = -1889.1072

$ws_ARM.Range("H61").Value = 4886.967
$ws_ARM.Range("I61").Value = 5174.8335
$ws_ARM.Range("J61").Value = 3735.5
$ws_ARM.Range("K61").Value = 5174.8335
$ws_ARM.Range("L61").Value = 3735.5
$ws_ARM.Range("M61").Value = -4962.8335
$ws_ARM.Range("N61").Value = -4159.5

$ws_ARM.Range("H74").Value = 1881.2727
$ws_ARM.Range("I74").Value = 1729.4
$ws_ARM.Range("J74").Value = 3400
$ws_ARM.Range("K74").Value = 1729.4
$ws_ARM.Range("L74").Value = 3400
$ws_ARM.Range("M74").Value = -855.4000000000001
$ws_ARM.Range("N74").Value = -5148

$ws_ARM.Range("H77").Value = 1881.2727
$ws_ARM.Range("I77").Value = 1729.4
$ws_ARM.Range("J77").Value = 3400
$ws_ARM.Range("K77").Value = 8647
$ws_ARM.Range("L77").Value = 17000
$ws_ARM.Range("M77").Value = -4279
$ws_ARM.Range("N77").Value = -25736

$ws_ARM.Range("H110").Value = 1694.9524
$ws_ARM.Range("I110").Value = 1723.5385
$ws_ARM.Range("K110").Value = 1723.5385
$ws_ARM.Range("M110").Value = 321.4614999999999

$ws_ARM.Range("H116").Value = 13726.182
$ws_ARM.Range("J116").Value = 3443.111
$ws_ARM.Range("L116").Value = 3443.111
$ws_ARM.Range("N116").Value = -8031.111

$ws_ARM.Range("H122").Value = 1082171.8
$ws_ARM.Range("I122").Value = 5374
$ws_ARM.Range("J122").Value = 1560748.5
$ws_ARM.Range("K122").Value = 16122
$ws_ARM.Range("L122").Value = 4682245.5
$ws_ARM.Range("M122").Value = -13672
$ws_ARM.Range("N122").Value = -4687145.5

$ws_ARM.Range("H136").Value = 4886.967
$ws_ARM.Range("I136").Value = 5174.8335
$ws_ARM.Range("J136").Value = 3735.5
$ws_ARM.Range("K136").Value = 15524.5005
$ws_ARM.Range("L136").Value = 11206.5
$ws_ARM.Range("M136").Value = -12974.5005
$ws_ARM.Range("N136").Value = -16306.5

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H3").Value = 13726.182
$ws_BSM.Range("J3").Value = 3443.111
$ws_BSM.Range("L3").Value = 3443.111
$ws_BSM.Range("N3").Value = -3671.111

$ws_BSM.Range("H20").Value = 1847
$ws_BSM.Range("I20").Value = 1124.05
$ws_BSM.Range("K20").Value = 1124.05
$ws_BSM.Range("M20").Value = -877.05

$ws_BSM.Range("H86").Value = 11442.333
$ws_BSM.Range("I86").Value = 14061.875
$ws_BSM.Range("J86").Value = 6203.25
$ws_BSM.Range("K86").Value = 14061.875
$ws_BSM.Range("L86").Value = 6203.25
$ws_BSM.Range("M86").Value = -12938.875
$ws_BSM.Range("N86").Value = -8449.25

$ws_BSM.Range("H89").Value = 11442.333
$ws_BSM.Range("I89").Value = 14061.875
$ws_BSM.Range("J89").Value = 6203.25
$ws_BSM.Range("K89").Value = 70309.375
$ws_BSM.Range("L89").Value = 31016.25
$ws_BSM.Range("M89").Value = -64693.375
$ws_BSM.Range("N89").Value = -42248.25

$ws_BSM.Range("H107").Value = 4190.8096
$ws_BSM.Range("I107").Value = 3767.5334
$ws_BSM.Range("K107").Value = 3767.5334
$ws_BSM.Range("M107").Value = -1847.5334

$ws_BSM.Range("H135").Value = 79355.8
$ws_BSM.Range("J135").Value = 79355.8
$ws_BSM.Range("L135").Value = 79355.8
$ws_BSM.Range("N135").Value = -89495.8

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 3874.7097
$ws_CRP.Range("I31").Value = 2724.08
$ws_CRP.Range("K31").Value = 2724.08
$ws_CRP.Range("M31").Value = -2429.08

$ws_CRP.Range("H34").Value = 3874.7097
$ws_CRP.Range("I34").Value = 2724.08
$ws_CRP.Range("K34").Value = 2724.08
$ws_CRP.Range("M34").Value = -2522.08

$ws_CRP.Range("H58").Value = 6890.9653
$ws_CRP.Range("I58").Value = 10525.786
$ws_CRP.Range("J58").Value = 3498.4666
$ws_CRP.Range("K58").Value = 10525.786
$ws_CRP.Range("L58").Value = 3498.4666
$ws_CRP.Range("M58").Value = -10322.786
$ws_CRP.Range("N58").Value = -3904.4666

$ws_CRP.Range("H94").Value = 2247.0588
$ws_CRP.Range("J94").Value = 1154.6923
$ws_CRP.Range("L94").Value = 1154.6923
$ws_CRP.Range("N94").Value = -2056.6923

$ws_CRP.Range("H98").Value = 34950
$ws_CRP.Range("J98").Value = 34950
$ws_CRP.Range("L98").Value = 34950
$ws_CRP.Range("N98").Value = -39442

$ws_CRP.Range("H108").Value = 31370.9
$ws_CRP.Range("J108").Value = 48796.8
$ws_CRP.Range("L108").Value = 48796.8
$ws_CRP.Range("N108").Value = -56476.8

$ws_CRP.Range("H136").Value = 6890.9653
$ws_CRP.Range("I136").Value = 10525.786
$ws_CRP.Range("J136").Value = 3498.4666
$ws_CRP.Range("K136").Value = 31577.358
$ws_CRP.Range("L136").Value = 10495.3998
$ws_CRP.Range("M136").Value = -29027.358
$ws_CRP.Range("N136").Value = -15595.3998

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H2").Value = 243.9375
$ws_CUL.Range("I2").Value = 194.57143
$ws_CUL.Range("K2").Value = 1167.42858
$ws_CUL.Range("M2").Value = -1054.42858

$ws_CUL.Range("H5").Value = 556933.9
$ws_CUL.Range("I5").Value = 1331.5385
$ws_CUL.Range("J5").Value = 2001500
$ws_CUL.Range("K5").Value = 3994.6155
$ws_CUL.Range("L5").Value = 6004500
$ws_CUL.Range("M5").Value = -3882.6155
$ws_CUL.Range("N5").Value = -6004724

$ws_CUL.Range("H68").Value = 2331.6667
$ws_CUL.Range("J68").Value = 2331.6667
$ws_CUL.Range("L68").Value = 6995.000100000001
$ws_CUL.Range("N68").Value = -8617.000100000001

$ws_CUL.Range("H71").Value = 2331.6667
$ws_CUL.Range("J71").Value = 2331.6667
$ws_CUL.Range("L71").Value = 20985.0003
$ws_CUL.Range("N71").Value = -29097.0003

$ws_CUL.Range("H125").Value = 19054.834
$ws_CUL.Range("I125").Value = 0
$ws_CUL.Range("K125").Value = 0
$ws_CUL.Range("M125").Value = ""

$ws_CUL.Range("H131").Value = 14926984
$ws_CUL.Range("I131").Value = 166668000
$ws_CUL.Range("J131").Value = 1639.2131
$ws_CUL.Range("K131").Value = 500004000
$ws_CUL.Range("L131").Value = 4917.6393
$ws_CUL.Range("M131").Value = -499998960
$ws_CUL.Range("N131").Value = -14997.6393

$ws_CUL.Range("H135").Value = 556933.9
$ws_CUL.Range("I135").Value = 1331.5385
$ws_CUL.Range("J135").Value = 2001500
$ws_CUL.Range("K135").Value = 11983.8465
$ws_CUL.Range("L135").Value = 18013500
$ws_CUL.Range("M135").Value = -9448.846500000001
$ws_CUL.Range("N135").Value = -18018570

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H80").Value = 2000.591
$ws_GSM.Range("I80").Value = 1804.6666
$ws_GSM.Range("J80").Value = 2420.4285
$ws_GSM.Range("K80").Value = 1804.6666
$ws_GSM.Range("L80").Value = 2420.4285
$ws_GSM.Range("M80").Value = -806.6666
$ws_GSM.Range("N80").Value = -4416.4285

$ws_GSM.Range("H83").Value = 2000.591
$ws_GSM.Range("I83").Value = 1804.6666
$ws_GSM.Range("J83").Value = 2420.4285
$ws_GSM.Range("K83").Value = 9023.333000000001
$ws_GSM.Range("L83").Value = 12102.1425
$ws_GSM.Range("M83").Value = -4031.333000000001
$ws_GSM.Range("N83").Value = -22086.1425

$ws_GSM.Range("H113").Value = 2346.2
$ws_GSM.Range("J113").Value = 2286.625
$ws_GSM.Range("L113").Value = 2286.625
$ws_GSM.Range("N113").Value = -6626.625

$ws_GSM.Range("H122").Value = 9260.723
$ws_GSM.Range("I122").Value = 10239.173
$ws_GSM.Range("J122").Value = 5207.143
$ws_GSM.Range("K122").Value = 30717.519
$ws_GSM.Range("L122").Value = 15621.429
$ws_GSM.Range("M122").Value = -28267.519
$ws_GSM.Range("N122").Value = -20521.429

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H22").Value = 3326.5386
$ws_LTW.Range("I22").Value = 2719.4285
$ws_LTW.Range("J22").Value = 4034.8333
$ws_LTW.Range("K22").Value = 2719.4285
$ws_LTW.Range("L22").Value = 4034.8333
$ws_LTW.Range("M22").Value = -2424.4285
$ws_LTW.Range("N22").Value = -4624.8333

$ws_LTW.Range("H27").Value = 3326.5386
$ws_LTW.Range("I27").Value = 2719.4285
$ws_LTW.Range("J27").Value = 4034.8333
$ws_LTW.Range("K27").Value = 2719.4285
$ws_LTW.Range("L27").Value = 4034.8333
$ws_LTW.Range("M27").Value = -2612.4285
$ws_LTW.Range("N27").Value = -4248.8333

$ws_LTW.Range("H100").Value = 6092.4
$ws_LTW.Range("I100").Value = 6853.273
$ws_LTW.Range("J100").Value = 4000
$ws_LTW.Range("K100").Value = 6853.273
$ws_LTW.Range("L100").Value = 4000
$ws_LTW.Range("M100").Value = -6312.273
$ws_LTW.Range("N100").Value = -5082

$ws_LTW.Range("H132").Value = 792822.1
$ws_LTW.Range("J132").Value = 5299.3335
$ws_LTW.Range("L132").Value = 15898.0005
$ws_LTW.Range("N132").Value = -20958.0005

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H40").Value = 54995
$ws_WVR.Range("I40").Value = 0
$ws_WVR.Range("J40").Value = 54995
$ws_WVR.Range("K40").Value = 0
$ws_WVR.Range("L40").Value = 54995
$ws_WVR.Range("M40").Value = ""
$ws_WVR.Range("N40").Value = -55293

$ws_WVR.Range("H62").Value = 85481.67999999999
$ws_WVR.Range("I62").Value = 163750.14
$ws_WVR.Range("J62").Value = 3299.8
$ws_WVR.Range("K62").Value = 163750.14
$ws_WVR.Range("L62").Value = 3299.8
$ws_WVR.Range("M62").Value = -163126.14
$ws_WVR.Range("N62").Value = -4547.8

$ws_WVR.Range("H65").Value = 85481.67999999999
$ws_WVR.Range("I65").Value = 163750.14
$ws_WVR.Range("J65").Value = 3299.8
$ws_WVR.Range("K65").Value = 818750.7000000001
$ws_WVR.Range("L65").Value = 16499
$ws_WVR.Range("M65").Value = -815630.7000000001
$ws_WVR.Range("N65").Value = -22739

$ws_WVR.Range("H81").Value = 29636.25
$ws_WVR.Range("I81").Value = 36279.25
$ws_WVR.Range("J81").Value = 16350.25
$ws_WVR.Range("K81").Value = 72558.5
$ws_WVR.Range("L81").Value = 32700.5
$ws_WVR.Range("M81").Value = -71497.5
$ws_WVR.Range("N81").Value = -34822.5

$ws_WVR.Range("H84").Value = 29636.25
$ws_WVR.Range("I84").Value = 36279.25
$ws_WVR.Range("J84").Value = 16350.25
$ws_WVR.Range("K84").Value = 362792.5
$ws_WVR.Range("L84").Value = 163502.5
$ws_WVR.Range("M84").Value = -357488.5
$ws_WVR.Range("N84").Value = -174110.5

$ws_WVR.Range("H122").Value = 2576.25
$ws_WVR.Range("I122").Value = 1890.8928
$ws_WVR.Range("K122").Value = 5672.678400000001
$ws_WVR.Range("M122").Value = -3222.678400000001

$ws_WVR.Range("H133").Value = 79999
$ws_WVR.Range("J133").Value = 79999
$ws_WVR.Range("L133").Value = 79999
$ws_WVR.Range("N133").Value = -90119

